$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.934.68"
$ws.Range("E2").Value = "  +0.50%  "
$ws.Range("D3").Value = "3.839.04"
$ws.Range("E3").Value = "  +0.92%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "688.54"
$ws.Range("E5").Value = "  +3.18%  "
$ws.Range("D6").Value = "172.71"
$ws.Range("E6").Value = "  +2.10%  "
$ws.Range("D7").Value = "3.836.42"
$ws.Range("E7").Value = "  +0.86%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "0.526"
$ws.Range("E9").Value = "  +0.08%  "
$ws.Range("E10").Value = "  +1.19%  "
$ws.Range("D11").Value = "7.38"
$ws.Range("E11").Value = "  +5.26%  "
$ws.Range("D12").Value = "0.460"
$ws.Range("E12").Value = "  -0.50%  "
$ws.Range("D13").Value = "0.0000258"
$ws.Range("E13").Value = "  +5.51%  "
$ws.Range("D14").Value = "36.54"
$ws.Range("E14").Value = "  +2.13%  "
$ws.Range("D15").Value = "4.499.79"
$ws.Range("E15").Value = "  +1.26%  "
$ws.Range("D16").Value = "3.845.93"
$ws.Range("D17").Value = "71.004.55"
$ws.Range("E17").Value = "  +0.69%  "
$ws.Range("D18").Value = "17.75"
$ws.Range("E18").Value = "  +0.22%  "
$ws.Range("D19").Value = "7.22"
$ws.Range("E19").Value = "  +0.39%  "
$ws.Range("E20").Value = "  +0.34%  "
$ws.Range("D21").Value = "11.13"
$ws.Range("E21").Value = "  -4.11%  "
$ws.Range("D22").Value = "486.62"
$ws.Range("E22").Value = "  +2.49%  "
$ws.Range("D23").Value = "0.720"
$ws.Range("E23").Value = "  +0.75%  "
$ws.Range("D24").Value = "84.60"
$ws.Range("E24").Value = "  +1.89%  "
$ws.Range("E25").Value = "  +2.05%  "
$ws.Range("D26").Value = "12.33"
$ws.Range("E26").Value = "  +1.02%  "
$ws.Range("D27").Value = "10.47"
$ws.Range("E27").Value = "  +1.03%  "
$ws.Range("B28").Value = "WrappedeETH"
$ws.Range("C28").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D28").Value = "4.002.94"
$ws.Range("E28").Value = "  +1.13%  "
$ws.Range("B29").Value = "Fetch.AI"
$ws.Range("C29").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D29").Value = "2.14"
$ws.Range("E29").Value = "  +0.47%  "
$ws.Range("E30").Value = "  +0.03%  "
$ws.Range("E31").Value = "  +8.77%  "
$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D32").Value = "7.62"
$ws.Range("E32").Value = "  +3.01%  "
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").Value = "2.31"
$ws.Range("E33").Value = "  +0.41%  "
$ws.Range("D34").Value = "29.69"
$ws.Range("E34").Value = "  +0.18%  "
$ws.Range("E35").Value = "  +3.31%  "
$ws.Range("D36").Value = "9.26"
$ws.Range("E36").Value = "  +1.47%  "
$ws.Range("D37").Value = "3.799.29"
$ws.Range("E37").Value = "  +1.01%  "
$ws.Range("D38").Value = "1.00"
$ws.Range("E38").Value = "  -0.06%  "
$ws.Range("E39").Value = "  +1.14%  "
$ws.Range("E40").Value = "  +12.66%  "
$ws.Range("D41").Value = "3.43"
$ws.Range("E41").Value = "  +0.25%  "
$ws.Range("D42").Value = "6.06"
$ws.Range("E42").Value = "  +1.53%  "
$ws.Range("E43").Value = "  +4.46%  "
$ws.Range("E44").Value = "  +0.05%  "
$ws.Range("D46").Value = "164.23"
$ws.Range("E46").Value = "  +3.51%  "
$ws.Range("E47").Value = "  +6.39%  "
$ws.Range("D48").Value = "48.65"
$ws.Range("E48").Value = "  +1.42%  "
$ws.Range("D49").Value = "44.39"
$ws.Range("E49").Value = "  -2.93%  "
$ws.Range("D50").Value = "0.302"
$ws.Range("E50").Value = "  +0.85%  "
$ws.Range("E51").Value = "  -3.06%  "
